$d = $word.ActiveDocument

$replacements = @(
    @("2024-06-14 Friday", "2024-06-15 Saturday"),
    @("149×6=", "474×5="),
    @("362×3=", "759×2="),
    @("561×7=", "311×8="),
    @("387×2=", "129×7="),
    @("678×9=", "608×6="),
    @("942×7=", "972×9="),
    @("432×6=", "699×5="),
    @("934×9=", "241×3="),
    @("251×5=", "765×2="),
    @("134×6=", "920×9="),
    @("695×8=", "442×3="),
    @("842×7=", "879×5="),
    @("230×2=", "733×2="),
    @("264×7=", "125×5="),
    @("243×8=", "280×5="),
    @("144×2=", "311×9="),
    @("625×6=", "119×4="),
    @("784×6=", "397×9="),
    @("157×5=", "780×9="),
    @("408×4=", "667×3="),
    @("406×8=", "201×6="),
    @("902×5=", "802×9="),
    @("517×3=", "558×6="),
    @("873×2=", "728×5="),
    @("387×6=", "225×6=")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
